$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 - match formatting of the other header cells
# (bold font, thin border all around, centered horizontally, top-aligned vertically)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# Data values for columns I and J, rows 2-48
$values = @{
    2  = @(7, 7)
    3  = @(7, 7)
    4  = @(6, 6)
    5  = @(5, 6)
    6  = @(7, 7)
    7  = @(3, 3)
    8  = @(8, 8)
    9  = @(7, 8)
    10 = @(8, 8)
    11 = @(9, 9)
    12 = @(6, 7)
    13 = @(6, 7)
    14 = @(6, 7)
    15 = @(7, 8)
    16 = @(7, 7)
    17 = @(5, 6)
    18 = @(9, 9)
    19 = @(7, 7)
    20 = @(9, 9)
    21 = @(5, 6)
    22 = @(6, 7)
    23 = @(7, 7)
    24 = @(7, 7)
    25 = @(5, 6)
    26 = @(6, 6)
    27 = @(5, 6)
    28 = @(7, 8)
    29 = @(8, 8)
    30 = @(8, 9)
    31 = @(6, 8)
    32 = @(10, 10)
    33 = @(6, 6)
    34 = @(7, 7)
    35 = @(8, 9)
    36 = @(5, 7)
    37 = @(6, 7)
    38 = @(5, 6)
    39 = @(8, 8)
    40 = @(8, 9)
    41 = @(6, 7)
    42 = @(4, 6)
    43 = @(1, 1)
    44 = @(6, 7)
    45 = @(8, 8)
    46 = @(4, 5)
    47 = @(4, 5)
    48 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
